# "todo listo antes de la presentacion"
# The underlying device time readings (time_inside / time_outside, in
# minutes) were recomputed into more precise decimal-hour figures. Update
# the source table (Sheet1!C2:D11) with the new values, refresh the
# PivotTable/PivotChart that summarizes them, tidy up the header font, and
# leave the selection where the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time_inside / time_outside values (decimal hours) for the 10 devices.
$newValues = @(
  @(178.6727533333333,  2.8337616666666672),
  @(201.55052444444439, 8.144715555555555),
  @(154.36643361111109, 3.3185705555555551),
  @(106.68880666666669, 17.479679444444439),
  @(73.316273888888887, 4.1962247222222224),
  @(117.35232999999999, 2.2162852777777782),
  @(42.655658611111107, 0.075203888888888876),
  @(183.83810361111111, 1.331275),
  @(212.19067861111111, 0.42446),
  @(127.88500500000001, 7.6480355555555546)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 3).Value = $newValues[$i][0]
  $ws.Cells.Item($row, 4).Value = $newValues[$i][1]
}

# Refresh the PivotTable (and its PivotChart) so the "Horas en casa" /
# "Horas afuera" summary picks up the new source data.
$pt = $ws.PivotTables("PivotTable1")
$pc = $pt.PivotCache()
$pc.Refresh()

# Header row (A1:D1) keeps its bold Calibri font, just normalized to carry
# an explicit font family id.
$ws.Range("A1:D1").Font.Family = 2

# Restore the author's last selection/view on the sheet.
$ws.Range("L8").Select()

Write-Output "time_in_and_out_home update applied"
